$d = $word.ActiveDocument

# "PENCACAHAN KEGIATAN LAPANGAN SURVEI BULAN FEBRUARI" -> split off "FEBRUARI"
# into its own run and highlight it yellow (matches the other ${...} template
# placeholders elsewhere in this document), leaving "BULAN " untouched.
$rng = $d.Content
$rng.Find.Execute("FEBRUARI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.HighlightColorIndex = 7

# Now that "FEBRUARI" is isolated in its own run, swap its text for the
# ${BULAN} placeholder token.
$rng2 = $d.Content
$rng2.Find.Execute("FEBRUARI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Text = "`${BULAN}"
